$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select the data range on sheet 1 before creating/populating the new sheet
$ws1.Range("B1:H33").Select() | Out-Null

# Add the new worksheet right after the existing one
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newWs.Name = "Blad1"

# Table 2 data (Instance, n, d, Apache total, SampLP total, SampLP Simplex, IterSampLP total, IterSampLP Simplex)
$data = New-Object 'object[,]' 33,8
$data[0,1] = 'n'
$data[0,2] = 'd'
$data[0,3] = 'Apache total'
$data[0,4] = ' SampLP total'
$data[0,5] = ' SampLP Simplex'
$data[0,6] = ' IterSampLP total'
$data[0,7] = ' IterSampLP Simplex'
$data[1,0] = 1
$data[1,1] = 1250
$data[1,2] = 5
$data[1,3] = 91605
$data[1,4] = 8571
$data[1,5] = 3663
$data[1,6] = 93396
$data[1,7] = 8672
$data[2,0] = 2
$data[2,1] = 1250
$data[2,2] = 5
$data[2,3] = 190610
$data[2,4] = 17437
$data[2,5] = 10941
$data[2,6] = 189220
$data[2,7] = 38456
$data[3,0] = 3
$data[3,1] = 1250
$data[3,2] = 5
$data[3,3] = 5599
$data[3,4] = 2080
$data[3,5] = 181
$data[3,6] = 17335
$data[3,7] = 234
$data[4,0] = 4
$data[4,1] = 1250
$data[4,2] = 5
$data[4,3] = 59826
$data[4,4] = 6676
$data[4,5] = 2950
$data[4,6] = 75248
$data[4,7] = 8101
$data[5,0] = 1
$data[5,1] = 2160
$data[5,2] = 6
$data[5,3] = 266887
$data[5,4] = 24334
$data[5,5] = 13178
$data[5,6] = 301377
$data[5,7] = 43118
$data[6,0] = 2
$data[6,1] = 2160
$data[6,2] = 6
$data[6,3] = 269059
$data[6,4] = 24017
$data[6,5] = 13366
$data[6,6] = 291398
$data[6,7] = 45011
$data[7,0] = 3
$data[7,1] = 2160
$data[7,2] = 6
$data[7,3] = 434566
$data[7,4] = 21049
$data[7,5] = 11186
$data[7,6] = 231637
$data[7,7] = 25728
$data[8,0] = 4
$data[8,1] = 2160
$data[8,2] = 6
$data[8,3] = 182747
$data[8,4] = 17205
$data[8,5] = 8382
$data[8,6] = 192605
$data[8,7] = 22020
$data[9,0] = 1
$data[9,1] = 270
$data[9,2] = 3
$data[9,3] = 3518
$data[9,4] = 843
$data[9,5] = 482
$data[9,6] = 6767
$data[9,7] = 1663
$data[10,0] = 2
$data[10,1] = 270
$data[10,2] = 3
$data[10,3] = 2498
$data[10,4] = 627
$data[10,5] = 314
$data[10,6] = 4419
$data[10,7] = 773
$data[11,0] = 3
$data[11,1] = 270
$data[11,2] = 3
$data[11,3] = 1249
$data[11,4] = 449
$data[11,5] = 189
$data[11,6] = 2677
$data[11,7] = 404
$data[12,0] = 4
$data[12,1] = 270
$data[12,2] = 3
$data[12,3] = 2591
$data[12,4] = 681
$data[12,5] = 353
$data[12,6] = 5213
$data[12,7] = 958
$data[13,0] = 1
$data[13,1] = 3430
$data[13,2] = 7
$data[13,3] = 1402486
$data[13,4] = 56464
$data[13,5] = 31761
$data[13,6] = 650137
$data[13,7] = 86166
$data[14,0] = 2
$data[14,1] = 3430
$data[14,2] = 7
$data[14,3] = 805804
$data[14,4] = 57138
$data[14,5] = 29486
$data[14,6] = 685525
$data[14,7] = 85931
$data[15,0] = 3
$data[15,1] = 3430
$data[15,2] = 7
$data[15,3] = 1209336
$data[15,4] = 55658
$data[15,5] = 27997
$data[15,6] = 753164
$data[15,7] = 86389
$data[16,0] = 4
$data[16,1] = 3430
$data[16,2] = 7
$data[16,3] = 1349811
$data[16,4] = 67305
$data[16,5] = 35871
$data[16,6] = 950145
$data[16,7] = 130263
$data[17,0] = 1
$data[17,1] = 5120
$data[17,2] = 8
$data[17,3] = 2087614
$data[17,4] = 86721
$data[17,5] = 32143
$data[17,6] = 1428837
$data[17,7] = 130015
$data[18,0] = 2
$data[18,1] = 5120
$data[18,2] = 8
$data[18,3] = 2891915
$data[18,4] = 87656
$data[18,5] = 36512
$data[18,6] = 1643975
$data[18,7] = 171009
$data[19,0] = 3
$data[19,1] = 5120
$data[19,2] = 8
$data[19,3] = 1170917
$data[19,4] = 71039
$data[19,5] = 28429
$data[19,6] = 797186
$data[19,7] = 68344
$data[20,0] = 4
$data[20,1] = 5120
$data[20,2] = 8
$data[20,3] = 1545349
$data[20,4] = 78305
$data[20,5] = 23893
$data[20,6] = 1470287
$data[20,7] = 101232
$data[21,0] = 1
$data[21,1] = 640
$data[21,2] = 4
$data[21,3] = 24974
$data[21,4] = 3044
$data[21,5] = 1815
$data[21,6] = 27322
$data[21,7] = 4874
$data[22,0] = 2
$data[22,1] = 640
$data[22,2] = 4
$data[22,3] = 23685
$data[22,4] = 2223
$data[22,5] = 1317
$data[22,6] = 19794
$data[22,7] = 3663
$data[23,0] = 3
$data[23,1] = 640
$data[23,2] = 4
$data[23,3] = 12578
$data[23,4] = 2169
$data[23,5] = 1250
$data[23,6] = 18999
$data[23,7] = 3526
$data[24,0] = 4
$data[24,1] = 640
$data[24,2] = 4
$data[24,3] = 17826
$data[24,4] = 2374
$data[24,5] = 1348
$data[24,6] = 25389
$data[24,7] = 4220
$data[25,0] = 1
$data[25,1] = 7290
$data[25,2] = 9
$data[25,3] = 8944001
$data[25,4] = 252987
$data[25,5] = 95365
$data[25,6] = 5591116
$data[25,7] = 751477
$data[26,0] = 2
$data[26,1] = 7290
$data[26,2] = 9
$data[26,3] = 3770595
$data[26,4] = 171645
$data[26,5] = 31653
$data[26,6] = 4177553
$data[26,7] = 267840
$data[27,0] = 3
$data[27,1] = 7290
$data[27,2] = 9
$data[27,3] = 2190897
$data[27,4] = 171699
$data[27,5] = 28134
$data[27,6] = 3349625
$data[27,7] = 210678
$data[28,0] = 4
$data[28,1] = 7290
$data[28,2] = 9
$data[28,3] = 4572966
$data[28,4] = 201903
$data[28,5] = 45554
$data[28,6] = 4859949
$data[28,7] = 348285
$data[29,0] = 1
$data[29,1] = 80
$data[29,2] = 2
$data[29,3] = 31
$data[29,4] = 25
$data[29,5] = 10
$data[29,6] = 127
$data[29,7] = 13
$data[30,0] = 2
$data[30,1] = 80
$data[30,2] = 2
$data[30,3] = 106
$data[30,4] = 57
$data[30,5] = 30
$data[30,6] = 255
$data[30,7] = 54
$data[31,0] = 3
$data[31,1] = 80
$data[31,2] = 2
$data[31,3] = 31
$data[31,4] = 25
$data[31,5] = 9
$data[31,6] = 135
$data[31,7] = 14
$data[32,0] = 4
$data[32,1] = 80
$data[32,2] = 2
$data[32,3] = 31
$data[32,4] = 25
$data[32,5] = 10
$data[32,6] = 128
$data[32,7] = 14

$newWs.Range("A1:H33").Value2 = $data

# Restore the selection on sheet 1 and set the new sheet's selection/active state
$ws1.Range("B1:H33").Select() | Out-Null
$newWs.Range("A34").Select() | Out-Null
